# "added Code and MenuDoc" -- add the encoder-input / INA-address lookup
# table (columns K:O) next to the existing A1/A0 address table.
#
# Cell-write order below deliberately mirrors the order the values were
# first typed by the original author so the shared-string table comes out
# in the same sequence (input, HEX add, 0x4A, 0x4F, testet?, v, 0x49, 0x4E,
# 0x46, 0x4B, 0x45, 0x48, Groen, orange).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row (row 3) -----------------------------------------------
$ws.Range("K3").Value = "input"
$ws.Range("L3").Value = "HEX add"

# --- first two HEX addresses (rows 4-5) --------------------------------
$ws.Range("L4").Value = "0x4A"
$ws.Range("L5").Value = "0x4F"

# --- remaining header cells (row 3) ------------------------------------
$ws.Range("M3").Value = "testet?"

# --- "tested" marks ------------------------------------------------------
$ws.Range("M4").Value = "v"

# --- rest of the HEX address column ------------------------------------
$ws.Range("L6").Value = "0x49"
$ws.Range("L7").Value = "0x4E"
$ws.Range("L8").Value = "0x46"
$ws.Range("L9").Value = "0x4B"
$ws.Range("L10").Value = "0x45"
$ws.Range("L11").Value = "0x48"

# --- colour legend (row 13) --------------------------------------------
$ws.Range("N13").Value = "Grøn"
$ws.Range("O13").Value = "orange"

# --- input index column (K4:K11) ---------------------------------------
$ws.Range("K4").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 4
$ws.Range("K8").Value = 5
$ws.Range("K9").Value = 6
$ws.Range("K10").Value = 7
$ws.Range("K11").Value = 8

# --- remaining "tested" marks -------------------------------------------
$ws.Range("M5").Value = "v"
$ws.Range("M6").Value = "v"
$ws.Range("M7").Value = "v"
$ws.Range("M8").Value = "v"
$ws.Range("M9").Value = "v"
$ws.Range("M11").Value = "v"

# --- address-pin colour columns (N/O), reusing A1/A0/Vs/SDA/SCL/GND ----
$ws.Range("N3").Value = "A1"
$ws.Range("O3").Value = "A0"

$ws.Range("N4").Value = "SDA"
$ws.Range("O4").Value = "SDA"

$ws.Range("N5").Value = "SCL"
$ws.Range("O5").Value = "SCL"

$ws.Range("N6").Value = "SDA"
$ws.Range("O6").Value = "Vs"

$ws.Range("N7").Value = "SCL"
$ws.Range("O7").Value = "SDA"

$ws.Range("N8").Value = "Vs"
$ws.Range("O8").Value = "SDA"

$ws.Range("N9").Value = "SDA"
$ws.Range("O9").Value = "SCL"

$ws.Range("N10").Value = "Vs"
$ws.Range("O10").Value = "Vs"

$ws.Range("N11").Value = "SDA"
$ws.Range("O11").Value = "GND"

# --- final selection, matching the author's last-saved cursor position -
$ws.Range("M10").Select()
